# Sentiment baselines were wrong because ru, tr and ja were included.
# Remove the Russian, Japanese and Turkish rows from every sheet (which also
# prunes those three entries out of the shared-string table on save), and
# refresh the remaining per-language baseline numbers.

$wb = $excel.ActiveWorkbook

# Final baseline values (after removing ru/ja/tr) for rows 4..13 on each sheet,
# keyed by the worksheet name.
$newValues = @{
    "Accuracy" = @{
        4  = 0.8768796992481203
        5  = 0.7803203661327232
        6  = 0.6045710139669871
        7  = 0.5138686131386861
        8  = 0.5921501706484642
        9  = 0.7455919395465995
        10 = 0.8458149779735683
        11 = 0.58207343412527
        12 = 0.8325724493846764
        13 = 0.6817817014446228
    }
    "Macro_Precision" = @{
        4  = 0.4384398496240601
        5  = 0.3901601830663616
        6  = 0.3022855069834935
        7  = 0.2569343065693431
        8  = 0.2960750853242321
        9  = 0.3727959697732998
        10 = 0.4229074889867842
        11 = 0.291036717062635
        12 = 0.4162862246923382
        13 = 0.3408908507223114
    }
    "Macro_Recall" = @{
        4  = 0.5
        5  = 0.5
        6  = 0.5
        7  = 0.5
        8  = 0.5
        9  = 0.5
        10 = 0.5
        11 = 0.5
        12 = 0.5
        13 = 0.5
    }
    "Macro_F1" = @{
        4  = 0.4672008012018027
        5  = 0.4383033419023136
        6  = 0.3767804657472304
        7  = 0.3394406943105111
        8  = 0.3719185423365488
        9  = 0.4271284271284271
        10 = 0.4582338902147972
        11 = 0.3679180887372014
        12 = 0.4543189818575684
        13 = 0.4053925077547125
    }
}

foreach ($sheetName in @("Accuracy", "Macro_Precision", "Macro_Recall", "Macro_F1")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 4 was "Russian" (deleting it shifts every later row up by one, so
    # the old "Japanese" row 12 becomes row 11, and old "Turkish" row 14
    # becomes row 13).
    $ws.Rows.Item(4).Delete() | Out-Null
    # Row 11 is now "Japanese"; deleting it shifts later rows up by one again,
    # so the old "Turkish" row (now 13) becomes row 12.
    $ws.Rows.Item(11).Delete() | Out-Null
    # Row 12 is now "Turkish".
    $ws.Rows.Item(12).Delete() | Out-Null

    $rowValues = $newValues[$sheetName]
    foreach ($r in $rowValues.Keys) {
        $ws.Cells.Item($r, 2).Value = $rowValues[$r]
    }
}
